$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '30.469.05'
$ws.Cells.Item(2, 5).Value = '  +0.12%  '

$ws.Cells.Item(3, 4).Value = '1.931.91'
$ws.Cells.Item(3, 5).Value = '  +4.34%  '

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = '@'
$c.Value = '0.9998'
$c.ClearFormats()

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '240.72'
$c.ClearFormats()
$ws.Cells.Item(5, 5).Value = '  +3.29%  '

$ws.Cells.Item(6, 5).Value = '  -0.03%  '

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = '@'
$c.Value = '0.4762'
$c.ClearFormats()
$ws.Cells.Item(7, 5).Value = '  +0.48%  '

$ws.Cells.Item(8, 2).Value = 'OKB'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = '@'
$c.Value = '44.68'
$c.ClearFormats()
$ws.Cells.Item(8, 5).Value = '  +3.19%  '

$ws.Cells.Item(9, 2).Value = 'Cardano'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = '@'
$c.Value = '0.2876'
$c.ClearFormats()
$ws.Cells.Item(9, 5).Value = '  +4.51%  '

$ws.Cells.Item(10, 2).Value = 'Dogecoin'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = '@'
$c.Value = '0.06637'
$c.ClearFormats()
$ws.Cells.Item(10, 5).Value = '  +4.66%  '

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = '@'
$c.Value = '108.00'
$c.ClearFormats()
$ws.Cells.Item(11, 5).Value = '  +27.51%  '

$ws.Cells.Item(12, 2).Value = 'Solana'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = '@'
$c.Value = '19.12'
$c.ClearFormats()
$ws.Cells.Item(12, 5).Value = '  +6.82%  '

$ws.Cells.Item(13, 2).Value = 'WrappedEther'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(13, 4).Value = '1.909.76'
$ws.Cells.Item(13, 5).Value = '  +3.11%  '

$ws.Cells.Item(14, 2).Value = 'TRON'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = '@'
$c.Value = '0.07612'
$c.ClearFormats()
$ws.Cells.Item(14, 5).Value = '  +1.92%  '

$ws.Cells.Item(15, 2).Value = 'Polkadot'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = '@'
$c.Value = '5.176'
$c.ClearFormats()
$ws.Cells.Item(15, 5).Value = '  +3.99%  '

$ws.Cells.Item(16, 2).Value = 'Polygon'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = '@'
$c.Value = '0.6614'
$c.ClearFormats()
$ws.Cells.Item(16, 5).Value = '  +6.17%  '

$ws.Cells.Item(17, 2).Value = 'BitcoinCash'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = '@'
$c.Value = '307.23'
$c.ClearFormats()
$ws.Cells.Item(17, 5).Value = '  +24.62%  '

$ws.Cells.Item(18, 2).Value = 'WrappedBTC'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(18, 4).Value = '30.486.39'
$ws.Cells.Item(18, 5).Value = '  +0.35%  '

$ws.Cells.Item(19, 2).Value = 'Avalanche'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = '@'
$c.Value = '13.01'
$c.ClearFormats()
$ws.Cells.Item(19, 5).Value = '  +2.57%  '

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = '@'
$c.Value = '1.000'
$c.ClearFormats()
$ws.Cells.Item(20, 5).Value = '  -0.01%  '

$ws.Cells.Item(21, 2).Value = 'ShibaInu'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = '0.000007585'
$c.ClearFormats()
$ws.Cells.Item(21, 5).Value = '  +3.40%  '

$ws.Cells.Item(22, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(22, 4).Value = '2.171.61'
$ws.Cells.Item(22, 5).Value = '  +3.70%  '

$ws.Cells.Item(23, 2).Value = 'Uniswap'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Value = '5.297'
$c.ClearFormats()
$ws.Cells.Item(23, 5).Value = '  +7.45%  '

$ws.Cells.Item(24, 2).Value = 'BinanceUSD'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.ClearFormats()
$ws.Cells.Item(24, 5).Value = '  -0.04%  '

$ws.Cells.Item(25, 2).Value = 'Chainlink'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = '@'
$c.Value = '6.307'
$c.ClearFormats()
$ws.Cells.Item(25, 5).Value = '  +6.98%  '

$ws.Cells.Item(26, 2).Value = 'Cosmos'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = '@'
$c.Value = '9.328'
$c.ClearFormats()
$ws.Cells.Item(26, 5).Value = '  +3.55%  '

$ws.Cells.Item(27, 2).Value = 'Monero'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = '@'
$c.Value = '167.86'
$c.ClearFormats()
$ws.Cells.Item(27, 5).Value = '  +2.35%  '

$ws.Cells.Item(28, 2).Value = 'EthereumClassic'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = '@'
$c.Value = '20.57'
$c.ClearFormats()
$ws.Cells.Item(28, 5).Value = '  +14.45%  '

$ws.Cells.Item(29, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = '@'
$c.Value = '2.053'
$c.ClearFormats()
$ws.Cells.Item(29, 5).Value = '  +9.54%  '

$ws.Cells.Item(30, 2).Value = 'Stellar'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = '@'
$c.Value = '0.1106'
$c.ClearFormats()
$ws.Cells.Item(30, 5).Value = '  +8.04%  '

$ws.Cells.Item(31, 2).Value = 'Toncoin'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = '@'
$c.Value = '1.370'
$c.ClearFormats()
$ws.Cells.Item(31, 5).Value = '  +1.77%  '

$ws.Cells.Item(32, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = '@'
$c.Value = '4.093'
$c.ClearFormats()
$ws.Cells.Item(32, 5).Value = '  +1.55%  '

$ws.Cells.Item(33, 2).Value = 'Filecoin'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Value = '3.943'
$c.ClearFormats()
$ws.Cells.Item(33, 5).Value = '  +2.95%  '

$ws.Cells.Item(34, 2).Value = 'Hedera'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = '@'
$c.Value = '0.05024'
$c.ClearFormats()
$ws.Cells.Item(34, 5).Value = '  +4.27%  '

$ws.Cells.Item(35, 2).Value = 'ImmutableX'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = '@'
$c.Value = '0.7429'
$c.ClearFormats()
$ws.Cells.Item(35, 5).Value = '  +6.67%  '

$ws.Cells.Item(36, 2).Value = 'ARBITRUM'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = '@'
$c.Value = '1.157'
$c.ClearFormats()
$ws.Cells.Item(36, 5).Value = '  +2.57%  '

$ws.Cells.Item(37, 2).Value = 'HuobiToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = '@'
$c.Value = '2.751'
$c.ClearFormats()
$ws.Cells.Item(37, 5).Value = '  +1.89%  '

$ws.Cells.Item(38, 2).Value = 'VeChain'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = '@'
$c.Value = '0.01964'
$c.ClearFormats()
$ws.Cells.Item(38, 5).Value = '  +3.78%  '

$ws.Cells.Item(39, 2).Value = 'MXToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = '@'
$c.Value = '2.689'
$c.ClearFormats()
$ws.Cells.Item(39, 5).Value = '  +0.45%  '

$ws.Cells.Item(40, 2).Value = 'RenderToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = '@'
$c.Value = '2.042'
$c.ClearFormats()
$ws.Cells.Item(40, 5).Value = '  +3.20%  '

$ws.Cells.Item(41, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = '@'
$c.Value = '0.8815'
$c.ClearFormats()
$ws.Cells.Item(41, 5).Value = '  +0.93%  '

$ws.Cells.Item(42, 2).Value = 'Quant'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = '@'
$c.Value = '107.79'
$c.ClearFormats()
$ws.Cells.Item(42, 5).Value = '  +1.51%  '

$ws.Cells.Item(43, 2).Value = 'Aave'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = '@'
$c.Value = '70.31'
$c.ClearFormats()
$ws.Cells.Item(43, 5).Value = '  +11.26%  '

$ws.Cells.Item(44, 2).Value = 'FraxShare'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = '@'
$c.Value = '5.787'
$c.ClearFormats()
$ws.Cells.Item(44, 5).Value = '  +5.28%  '

$ws.Cells.Item(45, 2).Value = 'PaxDollar'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Value = '0.9999'
$c.ClearFormats()
$ws.Cells.Item(45, 5).Value = '  -0.07%  '

$ws.Cells.Item(46, 2).Value = 'TheSandbox'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = '@'
$c.Value = '0.4201'
$c.ClearFormats()
$ws.Cells.Item(46, 5).Value = '  +3.56%  '

$ws.Cells.Item(47, 2).Value = 'Aptos'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = '@'
$c.Value = '7.279'
$c.ClearFormats()
$ws.Cells.Item(47, 5).Value = '  +1.56%  '

$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = '@'
$c.Value = '9.264'
$c.ClearFormats()
$ws.Cells.Item(48, 5).Value = '  +8.84%  '

$ws.Cells.Item(49, 2).Value = 'Algorand'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = '@'
$c.Value = '0.1214'
$c.ClearFormats()
$ws.Cells.Item(49, 5).Value = '  +1.41%  '

$ws.Cells.Item(50, 2).Value = 'Elrond'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = '@'
$c.Value = '34.87'
$c.ClearFormats()
$ws.Cells.Item(50, 5).Value = '  +2.48%  '

$ws.Cells.Item(51, 2).Value = 'Cronos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = '@'
$c.Value = '0.05625'
$c.ClearFormats()
$ws.Cells.Item(51, 5).Value = '  +2.31%  '
